# Append the new data row (2025/12/25) to the mods-count table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Figure out the next empty row right after the current data (row 44 -> 45).
$newRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1

# Date column: write as literal text (leading apostrophe forces text instead of
# Excel's automatic date conversion), matching the other rows in the sheet.
$ws.Cells.Item($newRow, 1).Value = "'2025/12/25"
$ws.Cells.Item($newRow, 2).Value = "逃离鸭科夫"
$ws.Cells.Item($newRow, 3).Value = 1100

# Match the centered alignment style used by the rest of the table rows.
$newRowRange = $ws.Range("A" + $newRow + ":C" + $newRow)
$newRowRange.HorizontalAlignment = -4108
$newRowRange.VerticalAlignment = -4108
